# Github Auto Build at 2023-12-11 10:22
# Append two new cost-log rows (178 and 179) to the sheet, following the
# existing pattern of column A = timestamp (text) and column B = cost (number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(178, 1).Value = "2023-12-11 10:22:13"
$ws.Cells.Item(178, 2).Value = 0.0012

$ws.Cells.Item(179, 1).Value = "2023-12-11 10:22:32"
$ws.Cells.Item(179, 2).Value = 0.001
